# "map string object correction"
# The jhi_user sheet gets a new data row (row 2) built from a user record
# map: only the fields present for this user are written as shared-string
# values, while fields the map had no entry for (city, state, phoneNumber,
# profiles) are still "touched" (present as blank cells) without a value.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("jhi_user")

# firstName, lastName, address
$ws.Cells.Item(2, 2).Value = "Era"
$ws.Cells.Item(2, 3).Value = "Rice"
$ws.Cells.Item(2, 4).Value = "526 Gary Cape"

# city / state -> no value in the map, but the cell is still touched (blank)
$ws.Cells.Item(2, 6).Font.Bold = $false
$ws.Cells.Item(2, 7).Font.Bold = $false

# mobilePhoneNumber
$ws.Cells.Item(2, 8).Value = "107-834-0930"

# phoneNumber -> no value in the map, blank cell
$ws.Cells.Item(2, 9).Font.Bold = $false

# email, password
$ws.Cells.Item(2, 11).Value = "omer.gutmann@yahoo.com"
$ws.Cells.Item(2, 12).Value = "nQ6,LjR"

# profiles -> no value in the map, blank cell
$ws.Cells.Item(2, 14).Font.Bold = $false
